$wb = $excel.ActiveWorkbook

$wsTypo = $wb.Worksheets.Item("Typography")
$wsTrans = $wb.Worksheets.Item("Translation")

# --- Typography sheet: row 6 updates ---
$wsTypo.Range("D6").Value = 15
$wsTypo.Range("I6").Value = "0x0020-0x007F,0x00C0-0x00FF,0x0018"

# --- Translation sheet: row 73 update ---
$wsTrans.Range("F73").Value = "alpha [" + [char]176 + "]:"

# --- Translation sheet: new rows 91-97 ---
$wsTrans.Range("B91").Value = "STATUSMSG_SETTINGS_UART_TX_ERR"
$wsTrans.Range("C91").Value = "Default"
$wsTrans.Range("D91").Value = "Left"
$wsTrans.Range("E91").Value = "LTR"
$wsTrans.Range("F91").Value = "Napaka pri UART TX komunikaciji."
$wsTrans.Range("G91").Value = "UART TX communication error."

$wsTrans.Range("B92").Value = "STATUSMSG_SETTINGS_UART_TX_NOT_OKED"
$wsTrans.Range("C92").Value = "Default"
$wsTrans.Range("D92").Value = "Left"
$wsTrans.Range("E92").Value = "LTR"
$wsTrans.Range("F92").Value = "GRBL je zavrnil G-komando."
$wsTrans.Range("G92").Value = "GRBL denied G-code."

$wsTrans.Range("B93").Value = "SingleUseId101"
$wsTrans.Range("C93").Value = "Default"
$wsTrans.Range("D93").Value = "Left"
$wsTrans.Range("E93").Value = "LTR"
$wsTrans.Range("F93").Value = "[rez / na 0]"
$wsTrans.Range("G93").Value = "[cut / to 0]"

$wsTrans.Range("B94").Value = "SingleUseId102"
$wsTrans.Range("C94").Value = "Small"
$wsTrans.Range("D94").Value = "Left"
$wsTrans.Range("E94").Value = "LTR"
$wsTrans.Range("F94").Value = " <uartConsoleBfr>"
$wsTrans.Range("G94").Value = " <uartConsoleBfr_5>"

$wsTrans.Range("B95").Value = "SingleUseId103"
$wsTrans.Range("C95").Value = "Small"
$wsTrans.Range("D95").Value = "Left"
$wsTrans.Range("E95").Value = "LTR"
$wsTrans.Range("F95").Value = " <uartConsoleBfr>"
$wsTrans.Range("G95").Value = " <uartConsoleBfr_5>"

$wsTrans.Range("B96").Value = "SingleUseId104"
$wsTrans.Range("C96").Value = "Small"
$wsTrans.Range("D96").Value = "Left"
$wsTrans.Range("E96").Value = "LTR"
$wsTrans.Range("F96").Value = " <uartConsoleBfr>"
$wsTrans.Range("G96").Value = " <uartConsoleBfr_5>"

$wsTrans.Range("B97").Value = "SingleUseId105"
$wsTrans.Range("C97").Value = "Small"
$wsTrans.Range("D97").Value = "Left"
$wsTrans.Range("E97").Value = "LTR"
$wsTrans.Range("F97").Value = " <uartConsoleBfr>"
$wsTrans.Range("G97").Value = " <uartConsoleBfr_5>"
